$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the "Application as {{JOB_TITLE}}" heading paragraph and replace
#    its content: drop the spell-check proofErr wrapper tags and merge the
#    "Application" / " " / "as" runs into a single "Application as" run.
# ---------------------------------------------------------------------------
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Application as*") {
        $headingIndex = $i
        break
    }
}
if ($headingIndex -eq -1) {
    throw "Could not find the 'Application as {{JOB_TITLE}}' paragraph"
}

$headingPara = $d.Paragraphs.Item($headingIndex)
$headingRange = $headingPara.Range

$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2DE8F180" w14:textId="4B40BFDE"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Application as</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> {{JOB_TITLE}}</w:t></w:r></w:p>'

$headingRange.InsertXML($headingXml)

# ---------------------------------------------------------------------------
# 2) Insert a brand-new paragraph right after the heading paragraph to hold
#    the cover-letter salutation placeholder.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item($headingIndex)
$headingPara.Range.InsertParagraphAfter()

$salutationPara = $d.Paragraphs.Item($headingIndex + 1)
$salutationRange = $salutationPara.Range

$salutationXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>{{COVER_LETTER_SALUTATION</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>}}</w:t></w:r></w:p>'

$salutationRange.InsertXML($salutationXml)

# ---------------------------------------------------------------------------
# 3) Locate the cover-letter body paragraph (contains {{COVER_LETTER_BODY}})
#    and replace its content: tag every run with lang="en-US" and append the
#    valediction placeholder runs after the existing trailing line breaks.
# ---------------------------------------------------------------------------
$bodyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*{{COVER_LETTER_BODY}}*") {
        $bodyIndex = $i
        break
    }
}
if ($bodyIndex -eq -1) {
    throw "Could not find the '{{COVER_LETTER_BODY}}' paragraph"
}

$bodyPara = $d.Paragraphs.Item($bodyIndex)
$bodyRange = $bodyPara.Range

$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="64757146" w14:textId="2D4CD8F9"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>{{COVER_LETTER_BODY}}</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>{{COVER_L</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>ETTER_</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="323B4C"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>VALEDICTION}}</w:t></w:r></w:p>'

$bodyRange.InsertXML($bodyXml)

Write-Output "Edit applied: heading=$headingIndex salutation=$($headingIndex + 1) body=$bodyIndex"
